$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "28_01_2024"
$ws.Range("F2").Value = 983
$ws.Range("F3").Value = 933
$ws.Range("F4").Value = 1108
$ws.Range("F5").Value = 2882

$ws.Range("F6").Font.Underline = $true
[void]$ws.Range("F6").Select()
